$d = $word.ActiveDocument

$replacements = @(
    @("2025-04-14 Monday", "2025-04-15 Tuesday"),
    @("869÷6=", "261÷5="),
    @("236÷6=", "770÷5="),
    @("355÷8=", "613÷6="),
    @("994÷2=", "137÷7="),
    @("533÷5=", "451÷7="),
    @("328÷7=", "805÷2="),
    @("195÷5=", "187÷8="),
    @("948÷9=", "393÷7="),
    @("920÷4=", "764÷3="),
    @("458÷2=", "831÷7="),
    @("452÷2=", "507÷4="),
    @("900÷2=", "805÷9="),
    @("429÷8=", "913÷3="),
    @("107÷3=", "313÷3="),
    @("216÷7=", "589÷2="),
    @("190÷5=", "789÷3="),
    @("410÷6=", "115÷7="),
    @("721÷7=", "100÷5="),
    @("762÷9=", "893÷8="),
    @("238÷5=", "349÷2="),
    @("679÷8=", "697÷6="),
    @("287÷5=", "826÷3="),
    @("752÷9=", "456÷6="),
    @("239÷8=", "341÷2="),
    @("443÷8=", "414÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
